# 16.1.3.xlsx — "Возраст"/"Образование" header row restructuring.
#
# The two summary rows under the "Items" header used bare "Age (in years)"
# and "Education" captions. They were reworded to "By age (in years)" /
# "By education" style captions (and their ru/ky equivalents), matching the
# phrasing used elsewhere in the sheet ("По территории", "Аймактар боюнча",
# ...). Updating the cell text lets Excel's shared-string table naturally
# drop the now-unused old strings and append the new ones.
#
# Order of assignment matters only for the position new strings land at in
# the shared-string table (cosmetic), so column C (the English caption) is
# set first for each row, then B (Russian), then A (Kyrgyz) — mirroring how
# the new strings were appended in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: "Age (in years)" / "Возраст (в годах)" -> "By age (in years)" / "По возрасту (в годах)"
$ws.Range("C18").Value = "By age (in years) "

# Row 28: "Education" / "Образование" -> "By education" / "По образованию"
$ws.Range("C28").Value = "By education"

$ws.Range("B18").Value = "По возрасту (в годах)"
$ws.Range("A18").Value = "Жаш курагы боюнча (жылдарда)"

$ws.Range("A28").Value = "Билими боюнча"
$ws.Range("B28").Value = "По образованию"

# The old selection (left over from an editing session on D4) is cleared;
# put the cursor back on A1, the sheet's natural resting position.
$ws.Range("A1").Select()
